$d = $word.ActiveDocument

# Helper-free approach: for each target phrase, locate it with Find, collapse the
# resulting range to its start point, then InsertBefore() the new lead-in text.
# InsertBefore() on a collapsed range creates a brand-new run ahead of the
# existing one, leaving the original run (and its formatting) untouched --
# exactly matching the diff's "add a sibling <w:r> before the existing <w:r>".

# 1) "Contatta l'assistenza" -> prefix with "FATTI: "
$r1 = $d.Content
$r1.Find.Execute("Contatta l’assistenza", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r1.Collapse(1)
$r1.InsertBefore("FATTI: ")

# 2) "Pubblica un nuovo itinerario in piattaforma" (first, non-registered variant) -> prefix with "FATTO: "
#    (Find matches the first occurrence in the document, which is this shorter variant,
#     since the "...registrato" paragraph comes later.)
$r2 = $d.Content
$r2.Find.Execute("Pubblica un nuovo itinerario in piattaforma", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r2.Collapse(1)
$r2.InsertBefore("FATTO: ")

# 3) "Pubblica un nuovo itinerario in piattaforma registrato" -> prefix with "FATTO: "
$r3 = $d.Content
$r3.Find.Execute("Pubblica un nuovo itinerario in piattaforma registrato", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r3.Collapse(1)
$r3.InsertBefore("FATTO: ")
